# Generate Report for Handoff
# Updates status from "In Translation" to "Ready for handoff" and refreshes
# the associated "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
# timestamps, widening the affected date/status columns to fit the new text.

$wb = $excel.ActiveWorkbook

# The new status text is wider than the old one, so the Status-ish columns
# are widened to fit. The engine snaps ColumnWidth to an internal 1/6-char
# pixel grid, so we request a value that lands on the closest grid point to
# the target width (~17.216 chars).
$targetColumnWidth = 16.333333333333332

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-19 20:44:32"
$wsOverview.Columns.Item(5).ColumnWidth = $targetColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $targetColumnWidth

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-19 20:44:28"
$wsZhCn.Columns.Item(3).ColumnWidth = $targetColumnWidth

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-19 20:44:32"
$wsDeDe.Columns.Item(3).ColumnWidth = $targetColumnWidth
